# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.538.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.882.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +7.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.881.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.396"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.49%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.412.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.439.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.60"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000190"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.881.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.031.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "511.62"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.75"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.58"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.55"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "185.11"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +8.76%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.08"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.69"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0928"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +9.73%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.694"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +17.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.582"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.76"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.79%  "
